{"js": "// Helper: replace a paragraph's run content with literal, distinct <w:r> runs\n// (built from `texts`, each becoming its own run) by round-tripping through\n// the flat-OPC \"insertOoxml\" gate. A plain `range.insertText()` call would\n// get coalesced into a single run by the host's adjacent-run merge, which is\n// exactly what we must avoid for the \"<image>foundation.svg\" style splits.\nfunction runsOoxml(texts) {\n  const body = texts\n    .map((t) => {\n      const esc = String(t)\n        .replace(/&/g, \"&amp;\")\n        .replace(/</g, \"&lt;\")\n        .replace(/>/g, \"&gt;\");\n      const preserve = /^\\s|\\s$/.test(t) ? ' xml:space=\"preserve\"' : \"\";\n      return `<w:r><w:t${preserve}>${esc}</w:t></w:r>`;\n    })\n    .join(\"\");\n  return (\n    '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    `<w:body><w:p>${body}</w:p></w:body>` +\n    \"</w:document></pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst p = paragraphs.items;\n\n// Work from the bottom of the document upward so inserting new paragraphs\n// never invalidates the indices of paragraphs still to be processed.\n\n// 7) \"Introduction to Linux\" (index 26) \u2014 mark where the page last broke.\np[26]\n  .getRange()\n  .insertOoxml(\n    '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n      '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n      \"<w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>Introduction to Linux</w:t></w:r></w:p></w:body>\" +\n      \"</w:document></pkg:xmlData></pkg:part></pkg:package>\",\n    \"Replace\"\n  );\n\n// 6) \"$\" + \"8\" + \"9.95 per month\" (index 21) -> single run.\np[21].getRange().insertText(\"$89.95 per month\", \"Replace\");\n\n// 5) \"Package#3\" (index 19) -> \"<image>professional.svg\" (6 runs) + new\n//    \"Alt: Faceless man in a suit\" paragraph right after it.\np[19].getRange().insertOoxml(runsOoxml([\"<im\", \"a\", \"g\", \"e\", \">\", \"professional.svg\"]), \"Replace\");\np[19].insertParagraph(\"Alt: Faceless man in a suit\", \"After\");\n\n// 4) \"Includes \" + \"the following courses:\" (index 13) -> single run.\np[13].getRange().insertText(\"Includes the following courses:\", \"Replace\");\n\n// 3) \"$\" + \"79\" + \".95 per month\" (index 11) -> single run.\np[11].getRange().insertText(\"$79.95 per month\", \"Replace\");\n\n// 2) \"Package #2\" (index 9) -> \"<image>engineer.svg\" (4 runs) + new\n//    \"Alt: nicely dressed man with a hard hat pointing at the ground\".\np[9].getRange().insertOoxml(runsOoxml([\"<im\", \"age\", \">\", \"engineer.svg\"]), \"Replace\");\np[9].insertParagraph(\"Alt: nicely dressed man with a hard hat pointing at the ground\", \"After\");\n\n// 1) \"Package #1\" (index 0) -> \"<image>foundation.svg\" (4 runs) + new\n//    \"Alt: image of a castle wall\".\np[0].getRange().insertOoxml(runsOoxml([\"<im\", \"age\", \">\", \"foundation.svg\"]), \"Replace\");\np[0].insertParagraph(\"Alt: image of a castle wall\", \"After\");\n\nawait context.sync();\n", "ps1": "# Word COM interop script.\n#\n# Splitting e.g. \"Package #1\" into the four literal runs\n# \"<im\" / \"age\" / \">\" / \"foundation.svg\" can't be done with plain\n# Range.Text / Range.InsertAfter calls: the host coalesces adjacent runs\n# that share identical formatting back into one run (the same behavior\n# real Word shows when you type into one spot in one editing session), so a\n# sequence of InsertAfter calls collapses right back into a single run.\n# Range.InsertXML (COM) / Range.insertOoxml (Office.js) accept a literal\n# WordOpenXML/flat-OPC payload and splice it in verbatim, which is the one\n# path that leaves the separate <w:r> runs alone.\nfunction New-RunsFlatOpc {\n    param([string[]]$Texts)\n\n    $runsXml = ($Texts | ForEach-Object {\n        $t = $_\n        $esc = $t -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'\n        if ($t -match '^\\s|\\s$') {\n            \"<w:r><w:t xml:space=`\"preserve`\">$esc</w:t></w:r>\"\n        } else {\n            \"<w:r><w:t>$esc</w:t></w:r>\"\n        }\n    }) -join ''\n\n    return '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        \"<w:body><w:p>$runsXml</w:p></w:body>\" +\n        '</w:document></pkg:xmlData></pkg:part></pkg:package>'\n}\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n\n# Walk from the bottom of the document upward so that inserting new\n# paragraphs never shifts the 1-based index of a paragraph still to come.\n\n# 7) \"Introduction to Linux\" (paragraph 27) -> stamp the render-break marker.\n$p27 = $paras.Item(27)\n$p27.Range.InsertXML('<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>Introduction to Linux</w:t></w:r></w:p></w:body>' +\n    '</w:document></pkg:xmlData></pkg:part></pkg:package>')\n\n# 6) \"$\" + \"8\" + \"9.95 per month\" (paragraph 22) -> one run.\n$p22 = $paras.Item(22)\n$p22.Range.Text = \"`$89.95 per month\"\n\n# 5) \"Package#3\" (paragraph 20) -> \"<image>professional.svg\" (6 runs) plus a\n#    new \"Alt: Faceless man in a suit\" paragraph right after it.\n$p20 = $paras.Item(20)\n$p20.Range.InsertXML((New-RunsFlatOpc @(\"<im\", \"a\", \"g\", \"e\", \">\", \"professional.svg\")))\n$p20 = $paras.Item(20)\n$p20.Range.InsertParagraphAfter()\n$paras.Item(21).Range.Text = \"Alt: Faceless man in a suit\"\n\n# 4) \"Includes \" + \"the following courses:\" (paragraph 14) -> one run.\n$p14 = $paras.Item(14)\n$p14.Range.Text = \"Includes the following courses:\"\n\n# 3) \"$\" + \"79\" + \".95 per month\" (paragraph 12) -> one run.\n$p12 = $paras.Item(12)\n$p12.Range.Text = \"`$79.95 per month\"\n\n# 2) \"Package #2\" (paragraph 10) -> \"<image>engineer.svg\" (4 runs) plus a new\n#    \"Alt: nicely dressed man with a hard hat pointing at the ground\".\n$p10 = $paras.Item(10)\n$p10.Range.InsertXML((New-RunsFlatOpc @(\"<im\", \"age\", \">\", \"engineer.svg\")))\n$p10 = $paras.Item(10)\n$p10.Range.InsertParagraphAfter()\n$paras.Item(11).Range.Text = \"Alt: nicely dressed man with a hard hat pointing at the ground\"\n\n# 1) \"Package #1\" (paragraph 1) -> \"<image>foundation.svg\" (4 runs) plus a new\n#    \"Alt: image of a castle wall\".\n$p1 = $paras.Item(1)\n$p1.Range.InsertXML((New-RunsFlatOpc @(\"<im\", \"age\", \">\", \"foundation.svg\")))\n$p1 = $paras.Item(1)\n$p1.Range.InsertParagraphAfter()\n$paras.Item(2).Range.Text = \"Alt: image of a castle wall\"\n"}
